$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28
$ws.Range("D28").Value = "NVIDIA Jeston 환경 셋팅 1-2편 (JetPack 설치 On Jeston Nano)"
$ws.Range("E28").Value = "https://ropiens.tistory.com/87"

# Row 36
$ws.Range("D36").Value = "Dive into BYOL"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/310"

# Row 41
$ws.Range("D41").Value = "Multi-datasource를 위한 SQL Engine"
$ws.Range("E41").Value = "http://cloudinsight.net/data/multi-datasource%eb%a5%bc-%ec%9c%84%ed%95%9c-sql-engine/"

# Row 51
$ws.Range("D51").Value = "블로거분들께 Deco 01 V2 타블렛 추천해봅니다"
$ws.Range("E51").Value = "https://bskyvision.com/1098"
